$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 3.286832544864788
$ws.Range("C2").Value = 117.745847958593
$ws.Range("D2").Value = 22.3905356188092
$ws.Range("E2").Value = 2195978.878461985
$ws.Range("G2").Value = 2196122.301678107

# Row 3
$ws.Range("B3").Value = 0.000001295275857016165
$ws.Range("C3").Value = 0.000000001689667739057654
$ws.Range("D3").Value = 3.537761648806719
$ws.Range("E3").Value = 1133.036916526867
$ws.Range("G3").Value = 1136.57467947264

# Row 4
$ws.Range("B4").Value = 3.286832544864788
$ws.Range("C4").Value = 1.655778082260271
$ws.Range("D4").Value = 0.7527432677738641
$ws.Range("E4").Value = 0.4942365360607697
$ws.Range("G4").Value = 6.189590430959694

# Row 5
$ws.Range("B5").Value = 1.455362044514542
$ws.Range("C5").Value = 10.34677158129881
$ws.Range("D5").Value = 0.7527432677738641
$ws.Range("E5").Value = 1133.036916526867
$ws.Range("G5").Value = 1145.591793420454
